# CryCompanywiseStockReport_1.xlsx - stock-count correction pass
#
# The report lists, per supplier/company, one row per item with:
#   A = Sl.No   B = Item Code   C = Item Name   D = Rate   E = MRP
#   F = Qty     G = Stock Value (= D * F)
# each company block ends with a "Sub Total:" row whose B is SUM(G) over the
# block, and the sheet ends with a "Sub Total:" / "Grand Total:" pair (B724 /
# B725) that is SUM(B) over every company sub-total.
#
# This pass corrects the counted quantity (column F) for a batch of items,
# recalculates the resulting stock value (column G = D * F), re-swaps two
# mis-ordered duplicate-batch row pairs (their B/E/F/G belong to each other),
# and rolls the corrections up through the Sub Total / Grand Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value2 = $value
}

# ---------------------------------------------------------------------
# 1) Per-item quantity corrections: new Qty (F), recalculated Value (G)
# ---------------------------------------------------------------------
$qtyFixes = @(
    @{ Row = 31;  F = 36;  G = 1844.28 }
    @{ Row = 67;  F = 198; G = 51624.54 }
    @{ Row = 158; F = 130; G = 6202.3 }
    @{ Row = 192; F = 11;  G = 476.08 }
    @{ Row = 212; F = 51;  G = 3304.8 }
    @{ Row = 221; F = 159; G = 17854.11 }
    @{ Row = 228; F = 486; G = 8991 }
    @{ Row = 285; F = 7;   G = 535.29 }
    @{ Row = 290; F = 141; G = 6608.67 }
    @{ Row = 331; F = 6;   G = 907.02 }
    @{ Row = 349; F = 49;  G = 9488.360000000001 }
    @{ Row = 393; F = 385; G = 37191 }
    @{ Row = 404; F = 83;  G = 4658.79 }
    @{ Row = 408; F = 24;  G = 823.4400000000001 }
    @{ Row = 414; F = 175; G = 2773.75 }
    @{ Row = 422; F = 51;  G = 1497.36 }
    @{ Row = 461; F = 34;  G = 7555.82 }
    @{ Row = 517; F = 199; G = 19874.13 }
    @{ Row = 518; F = 13;  G = 1541.54 }
    @{ Row = 525; F = 347; G = 19043.36 }
    @{ Row = 534; F = 241; G = 3822.26 }
    @{ Row = 535; F = 108; G = 3575.88 }
    @{ Row = 538; F = 4;   G = 172.72 }
    @{ Row = 564; F = 151; G = 18399.35 }
    @{ Row = 618; F = 222; G = 33391.02 }
    @{ Row = 665; F = 30;  G = 1606.2 }
    @{ Row = 680; F = 539; G = 87916.28999999999 }
)

foreach ($fix in $qtyFixes) {
    Set-Cell ("F{0}" -f $fix.Row) $fix.F
    Set-Cell ("G{0}" -f $fix.Row) $fix.G
}

# ---------------------------------------------------------------------
# 2) Duplicate-batch rows that were listed in the wrong order: the two
#    rows share the same item (A/C/D untouched) but their Code/MRP/Qty/
#    Value (B/E/F/G) belong to the other row - swap them back.
# ---------------------------------------------------------------------
function Swap-Rows($row1, $row2) {
    foreach ($col in @('B', 'E', 'F', 'G')) {
        $addr1 = "{0}{1}" -f $col, $row1
        $addr2 = "{0}{1}" -f $col, $row2
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        Set-Cell $addr1 $v2
        Set-Cell $addr2 $v1
    }
}

Swap-Rows 303 304
Swap-Rows 310 311
Swap-Rows 485 486
Swap-Rows 502 503

# ---------------------------------------------------------------------
# 3) Roll the corrections up into each affected company's "Sub Total:"
#    row (B = SUM of that company's G column).
# ---------------------------------------------------------------------
$subtotals = @(
    @{ Start = 24;  End = 39;  Row = 40 }
    @{ Start = 42;  End = 71;  Row = 72 }
    @{ Start = 155; End = 160; Row = 161 }
    @{ Start = 178; End = 198; Row = 199 }
    @{ Start = 212; End = 213; Row = 214 }
    @{ Start = 216; End = 223; Row = 224 }
    @{ Start = 226; End = 234; Row = 235 }
    @{ Start = 256; End = 300; Row = 301 }
    @{ Start = 303; End = 333; Row = 334 }
    @{ Start = 348; End = 361; Row = 362 }
    @{ Start = 386; End = 394; Row = 395 }
    @{ Start = 397; End = 422; Row = 423 }
    @{ Start = 456; End = 463; Row = 464 }
    @{ Start = 516; End = 530; Row = 531 }
    @{ Start = 533; End = 540; Row = 541 }
    @{ Start = 564; End = 566; Row = 567 }
    @{ Start = 615; End = 633; Row = 634 }
    @{ Start = 665; End = 673; Row = 674 }
    @{ Start = 679; End = 685; Row = 686 }
)

foreach ($st in $subtotals) {
    $rng = $ws.Range(("G{0}:G{1}" -f $st.Start, $st.End))
    $sectionSum = $excel.WorksheetFunction.Sum($rng)
    Set-Cell ("B{0}" -f $st.Row) $sectionSum
}

# ---------------------------------------------------------------------
# 4) Roll everything up into the sheet-level "Sub Total:" / "Grand
#    Total:" rows (724 / 725 = SUM of every company Sub Total's B).
# ---------------------------------------------------------------------
$allSubtotalRows = @(9,15,18,22,40,72,75,89,93,96,103,129,134,139,143,153,161,170,176,199,203,207,210,214,224,235,246,254,301,334,338,342,346,362,369,378,381,384,395,423,433,444,447,451,454,464,469,482,488,499,514,531,541,544,547,562,567,579,588,591,594,601,609,613,634,641,663,674,677,686,690,697,719,723)

$grandTotal = 0
foreach ($r in $allSubtotalRows) {
    $grandTotal += $ws.Range(("B{0}" -f $r)).Value2
}

Set-Cell "B724" $grandTotal
Set-Cell "B725" $grandTotal
